# Updates after-sales.xlsx: removes the superseded HBAR lot that used to sit
# at row 14 (qty 53.64615385 / basis 3.17) and refreshes every remaining sale
# with its precise timestamp (date + time of day) instead of a date-only value.
# A few quantity/cost-basis figures were also corrected at the same time.
# (commit: "added LIFO, FIFO, HIFO, and method for each year")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the obsolete row; Excel shifts rows 15-46 up to fill the gap,
# which also fixes the sheet dimension from A1:D46 down to A1:D45.
$ws.Rows.Item(14).Delete()

# Data for the shifted/updated rows, as (row, date-serial, currency, quantity, cost basis)
$rowData = @(
    ,@(2, 45345.19027777778, "XRP", 3616.65739, 1970.64)
    ,@(3, 45327.97185185185, "HBAR", 10.87677555, 0.74)
    ,@(4, 45320.48010416667, "HBAR", 72.68400164000001, 5.33)
    ,@(5, 45281.83146990741, "HBAR", 34.12963782, 2.97)
    ,@(6, 45266.91903935185, "HBAR", 48.30298109, 3.18)
    ,@(7, 45247.9331712963, "HBAR", 10.96719865, 0.67)
    ,@(8, 45243.92203703704, "HBAR", 55.78260692, 3.47)
    ,@(9, 45230.47232638889, "USDC", 0.159954, 0.16)
    ,@(10, 45224.69118055556, "HBAR", 15.68345492, 0.8100000000000001)
    ,@(11, 45219.75883101852, "HBAR", 297.4689965, 13.63)
    ,@(12, 45146.03603009259, "HBAR", 470.261323, 26.55)
    ,@(13, 45146.03285879629, "HBAR", 43801.36281517, 2561.94)
    ,@(14, 45131.44159722222, "XRP", 1418.744626, 1000)
    ,@(15, 45124.80702546296, "XRP", 922.768313, 674.8200000000001)
    ,@(16, 45124.77539351852, "XRP", 941.043309, 681.5)
    ,@(17, 45124.75659722222, "XRP", 971.856015, 725.9400000000001)
    ,@(18, 45087.11783564815, "HBAR", 130.2443386, 6.31)
    ,@(19, 45067.3153125, "HBAR", 58.46972882, 3.1)
    ,@(20, 45058.33905092593, "HBAR", 283.8796531, 14.59)
    ,@(21, 45058.19652777778, "XRP", 44.590986, 19.13)
    ,@(22, 45054.77490740741, "XRP", 113.287, 51.02)
    ,@(23, 45014.81280092592, "HBAR", 6.45060194, 0.38)
    ,@(24, 45013.15342592593, "XRP", 2066.699516, 992.3)
    ,@(25, 45012.1265162037, "XRP", 805.484107, 361.53)
    ,@(26, 45010.50967592592, "XRP", 219.938426, 100)
    ,@(27, 45006.65659722222, "XRP", 1152.377951, 503.82)
    ,@(28, 45005.71923611111, "HBAR", 317.43499761, 20)
    ,@(29, 44998.83645833333, "HBAR", 309.67942733, 20)
    ,@(30, 44995.76827546296, "HBAR", 10127.87723785, 600)
    ,@(31, 44991.70733796297, "HBAR", 311.96112243, 20)
    ,@(32, 44987.0953125, "HBAR", 12639.46303323, 872.4299999999999)
    ,@(33, 44984.69497685185, "HBAR", 274.33140336, 20)
    ,@(34, 44963.71480324074, "HBAR", 264.01674599, 18.51)
    ,@(35, 44956.84309027778, "HBAR", 288.44870791, 18.51)
    ,@(36, 44949.78057870371, "HBAR", 265.38112392, 18.51)
    ,@(37, 44946.91366898148, "HBAR", 6746.06763447, 513.26)
    ,@(38, 44942.81034722222, "HBAR", 372.58940273, 20)
    ,@(39, 44937.76524305555, "HBAR", 872.76773121, 38.93)
    ,@(40, 44937.76501157408, "HBAR", 655.50684043, 29.23)
    ,@(41, 44937.76446759259, "HBAR", 656.21331996, 29.24)
    ,@(42, 44935.6687037037, "HBAR", 443.1084391, 20)
    ,@(43, 44893.29744212963, "ETH", 0.0003157799999999981, 0.3499999999999979)
    ,@(44, 44867.11673611111, "HBAR", 839.14913695, 48.73)
    ,@(45, 44867.11646990741, "HBAR", 839.61242671, 48.74)
)

foreach ($row in $rowData) {
    $r = $row[0]
    $ws.Range("A$r").Value = $row[1]
    $ws.Range("B$r").Value = $row[2]
    $ws.Range("C$r").Value = $row[3]
    $ws.Range("D$r").Value = $row[4]
}
